$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginData")
$ws2 = $wb.Worksheets.Item("SearchData")

# LoginData: rename the "profileName" key to "osandaProfileName"
$ws1.Range("A4").Value = "osandaProfileName"

# SearchData: replace the search-question sample row with T-Shirt dress search data
$ws2.Range("A2").Value = "tShirtDress"
$ws2.Range("B2").Value = "T-Shirt"
$ws2.Range("A3").Value = ""
$ws2.Range("B3").Value = ""

# Update selections: LoginData is no longer the active tab, SearchData is now active
$ws1.Range("I14").Select()
$ws2.Activate()
$ws2.Range("G11").Select()
